$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12926
$ws1.Range("F5").Value = 89
$ws1.Range("F6").Value = 77
$ws1.Range("F8").Value = 22
$ws1.Range("F10").Value = 12880
$ws1.Range("F13").Value = 8694
$ws1.Range("F14").Value = 7696
$ws1.Range("F18").Value = 128
$ws1.Range("F19").Value = 984
$ws1.Range("F25").Value = 87

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12926
$ws4.Range("F6").Value = 89
$ws4.Range("F7").Value = 77
$ws4.Range("F9").Value = 22
$ws4.Range("F11").Value = 12880
$ws4.Range("F14").Value = 8694
$ws4.Range("F15").Value = 7696
$ws4.Range("F19").Value = 128
$ws4.Range("F20").Value = 984
$ws4.Range("F27").Value = 87
